$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 7482
$ws.Range("F4").Value = 3558
$ws.Range("F6").Value = 3878
$ws.Range("F9").Value = 81
$ws.Range("F10").Value = 114
$ws.Range("F11").Value = 167
$ws.Range("F12").Value = 519
$ws.Range("F13").Value = 20
$ws.Range("F18").Value = 4208
$ws.Range("F20").Value = 418
$ws.Range("F22").Value = 541
$ws.Range("F23").Value = 1911
$ws.Range("F24").Value = 119
$ws.Range("F25").Value = 104
$ws.Range("F26").Value = 70
$ws.Range("F27").Value = 3091
$ws.Range("F28").Value = 2342
$ws.Range("F31").Value = 96
$ws.Range("F32").Value = 102
$ws.Range("F33").Value = 127
$ws.Range("F37").Value = 4431
$ws.Range("F38").Value = 513
$ws.Range("F42").Value = 842
$ws.Range("F43").Value = 241
$ws.Range("F45").Value = 1676
$ws.Range("F46").Value = 266
$ws.Range("F48").Value = 619
$ws.Range("F49").Value = 730

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4
$ws.Range("F6").Value = 67
$ws.Range("F10").Value = 47
$ws.Range("F14").Value = 109
$ws.Range("F19").Value = 612

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 7482
$ws.Range("F6").Value = 3558
$ws.Range("F7").Value = 3878
$ws.Range("F10").Value = 81
$ws.Range("F11").Value = 114
$ws.Range("F13").Value = 167
$ws.Range("F14").Value = 519
$ws.Range("F15").Value = 67
$ws.Range("F20").Value = 4208
$ws.Range("F23").Value = 47
$ws.Range("F25").Value = 418
$ws.Range("F26").Value = 541
$ws.Range("F27").Value = 1911
$ws.Range("F28").Value = 119
$ws.Range("F29").Value = 104
$ws.Range("F30").Value = 3091
$ws.Range("F31").Value = 2342
$ws.Range("F34").Value = 127
$ws.Range("F36").Value = 109
$ws.Range("F37").Value = 4431
$ws.Range("F39").Value = 513
$ws.Range("F42").Value = 842
$ws.Range("F43").Value = 241
$ws.Range("F45").Value = 1676
$ws.Range("F46").Value = 266
$ws.Range("F48").Value = 619
$ws.Range("F49").Value = 730
